# Apply updates for "Add data for 2022-07-28":
# - Rename sheet from "Through 2022-07-19" to "Through 2022-07-20"
# - Update header label in I1 from "2022 (through 07-19)" to "2022 (through 07-20)"
# - Update I8 from 105 to 112
# - Update I14 (Total) from 911 to 918

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2022-07-20"

$ws.Range("I1").Value = "2022 (through 07-20)"
$ws.Range("I8").Value = 112
$ws.Range("I14").Value = 918
